$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 3573.238
$ws.Range("I88").Value = 1623
$ws.Range("J88").Value = 3778.5264
$ws.Range("K88").Value = 1623
$ws.Range("L88").Value = 3778.5264
$ws.Range("M88").Value = -1217
$ws.Range("N88").Value = -4590.526400000001
$ws.Range("H91").Value = 3573.238
$ws.Range("I91").Value = 1623
$ws.Range("J91").Value = 3778.5264
$ws.Range("K91").Value = 1623
$ws.Range("L91").Value = 3778.5264
$ws.Range("M91").Value = -219
$ws.Range("N91").Value = -6586.526400000001
$ws.Range("H98").Value = 1775.72
$ws.Range("I98").Value = 1475.6957
$ws.Range("K98").Value = 1475.6957
$ws.Range("M98").Value = 22.30430000000001
$ws.Range("H106").Value = 3609
$ws.Range("I106").Value = 3609
$ws.Range("K106").Value = 3609
$ws.Range("M106").Value = -2978
$ws.Range("H122").Value = 1775.72
$ws.Range("I122").Value = 1475.6957
$ws.Range("K122").Value = 4427.0871
$ws.Range("M122").Value = -1977.0871
$ws.Range("H132").Value = 25003050
$ws.Range("I132").Value = 27781022
$ws.Range("K132").Value = 83343066
$ws.Range("M132").Value = -83340536
$ws.Range("H135").Value = 1540.75
$ws.Range("I135").Value = 1186.375
$ws.Range("K135").Value = 10677.375
$ws.Range("M135").Value = -8142.375
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1846
$ws.Range("I63").Value = 1846
$ws.Range("K63").Value = 1846
$ws.Range("M63").Value = -1160
$ws.Range("H66").Value = 1846
$ws.Range("I66").Value = 1846
$ws.Range("K66").Value = 9230
$ws.Range("M66").Value = -5798
$ws.Range("H97").Value = 1474965.6
$ws.Range("I97").Value = 1707645
$ws.Range("K97").Value = 1707645
$ws.Range("M97").Value = -1707149
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H132").Value = 3737.3572
$ws.Range("I132").Value = 3372.5
$ws.Range("K132").Value = 10117.5
$ws.Range("M132").Value = -7587.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1279.3077
$ws.Range("I20").Value = 958.06665
$ws.Range("J20").Value = 1717.3636
$ws.Range("K20").Value = 958.06665
$ws.Range("L20").Value = 1717.3636
$ws.Range("M20").Value = -711.06665
$ws.Range("N20").Value = -2211.3636
$ws.Range("H86").Value = 3130203.8
$ws.Range("I86").Value = 3577140
$ws.Range("J86").Value = 1651.25
$ws.Range("K86").Value = 3577140
$ws.Range("L86").Value = 1651.25
$ws.Range("M86").Value = -3576017
$ws.Range("N86").Value = -3897.25
$ws.Range("H89").Value = 3130203.8
$ws.Range("I89").Value = 3577140
$ws.Range("J89").Value = 1651.25
$ws.Range("K89").Value = 17885700
$ws.Range("L89").Value = 8256.25
$ws.Range("M89").Value = -17880084
$ws.Range("N89").Value = -19488.25
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H134").Value = 4982.2085
$ws.Range("I134").Value = 2356.2307
$ws.Range("K134").Value = 7068.6921
$ws.Range("M134").Value = -4533.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1996.0769
$ws.Range("I16").Value = 1584.8889
$ws.Range("K16").Value = 1584.8889
$ws.Range("M16").Value = -1297.8889
$ws.Range("H86").Value = 9446.833000000001
$ws.Range("I86").Value = 8145.125
$ws.Range("J86").Value = 10488.2
$ws.Range("K86").Value = 8145.125
$ws.Range("L86").Value = 10488.2
$ws.Range("M86").Value = -7022.125
$ws.Range("N86").Value = -12734.2
$ws.Range("H89").Value = 9446.833000000001
$ws.Range("I89").Value = 8145.125
$ws.Range("J89").Value = 10488.2
$ws.Range("K89").Value = 40725.625
$ws.Range("L89").Value = 52441
$ws.Range("M89").Value = -35109.625
$ws.Range("N89").Value = -63673
$ws.Range("H105").Value = 1593.6522
$ws.Range("I105").Value = 1599.125
$ws.Range("J105").Value = 1581.1428
$ws.Range("K105").Value = 1599.125
$ws.Range("L105").Value = 1581.1428
$ws.Range("M105").Value = 147.875
$ws.Range("N105").Value = -5075.1428
$ws.Range("H113").Value = 1996.0769
$ws.Range("I113").Value = 1584.8889
$ws.Range("K113").Value = 1584.8889
$ws.Range("M113").Value = 585.1111000000001
$ws.Range("H134").Value = 37749.19
$ws.Range("I134").Value = 62519.332
$ws.Range("J134").Value = 3971.7273
$ws.Range("K134").Value = 187557.996
$ws.Range("L134").Value = 11915.1819
$ws.Range("M134").Value = -185022.996
$ws.Range("N134").Value = -16985.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6409.091
$ws.Range("I56").Value = 6409.091
$ws.Range("K56").Value = 6409.091
$ws.Range("M56").Value = -5879.091
$ws.Range("H107").Value = 327.4
$ws.Range("I107").Value = 173.41176
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 520.23528
$ws.Range("L107").Value = 3600
$ws.Range("M107").Value = 1399.76472
$ws.Range("N107").Value = -7440
$ws.Range("H141").Value = 2562.25
$ws.Range("I141").Value = 1999
$ws.Range("K141").Value = 5997
$ws.Range("M141").Value = -817

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 45000
$ws.Range("J34").Value = 45000
$ws.Range("L34").Value = 45000
$ws.Range("N34").Value = -45536
$ws.Range("H76").Value = 45000
$ws.Range("J76").Value = 45000
$ws.Range("L76").Value = 45000
$ws.Range("N76").Value = -45630
$ws.Range("H79").Value = 45000
$ws.Range("J79").Value = 45000
$ws.Range("L79").Value = 45000
$ws.Range("N79").Value = -47184
$ws.Range("H102").Value = 7998976.5
$ws.Range("I102").Value = 12348677
$ws.Range("K102").Value = 12348677
$ws.Range("M102").Value = -12347055
$ws.Range("H126").Value = 8052823
$ws.Range("I126").Value = 5053537.5
$ws.Range("J126").Value = 11909047
$ws.Range("K126").Value = 15160612.5
$ws.Range("L126").Value = 35727141
$ws.Range("M126").Value = -15158142.5
$ws.Range("N126").Value = -35732081
$ws.Range("H132").Value = 4223.8667
$ws.Range("J132").Value = 6247.8
$ws.Range("L132").Value = 18743.4
$ws.Range("N132").Value = -23803.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4791.7334
$ws.Range("I46").Value = 3084.1428
$ws.Range("K46").Value = 3084.1428
$ws.Range("M46").Value = -2896.1428
$ws.Range("H55").Value = 2206.2593
$ws.Range("I55").Value = 1900.9375
$ws.Range("J55").Value = 2650.3635
$ws.Range("K55").Value = 1900.9375
$ws.Range("L55").Value = 2650.3635
$ws.Range("M55").Value = -1727.9375
$ws.Range("N55").Value = -2996.3635
$ws.Range("H61").Value = 3369199.8
$ws.Range("I61").Value = 4117550.5
$ws.Range("K61").Value = 4117550.5
$ws.Range("M61").Value = -4117348.5
$ws.Range("H68").Value = 3733.3333
$ws.Range("I68").Value = 2750
$ws.Range("K68").Value = 2750
$ws.Range("M68").Value = -2001
$ws.Range("H71").Value = 3733.3333
$ws.Range("I71").Value = 2750
$ws.Range("K71").Value = 13750
$ws.Range("M71").Value = -10006
$ws.Range("H93").Value = 37040132
$ws.Range("I93").Value = 83337336
$ws.Range("J93").Value = 2369.6
$ws.Range("K93").Value = 83337336
$ws.Range("L93").Value = 2369.6
$ws.Range("M93").Value = -83336088
$ws.Range("N93").Value = -4865.6
$ws.Range("H113").Value = 3369199.8
$ws.Range("I113").Value = 4117550.5
$ws.Range("K113").Value = 4117550.5
$ws.Range("M113").Value = -4115380.5
$ws.Range("H122").Value = 5323.0435
$ws.Range("I122").Value = 3604
$ws.Range("J122").Value = 7557.8
$ws.Range("K122").Value = 10812
$ws.Range("L122").Value = 22673.4
$ws.Range("M122").Value = -8362
$ws.Range("N122").Value = -27573.4
$ws.Range("H136").Value = 120115.47
$ws.Range("I136").Value = 127547.69
$ws.Range("K136").Value = 382643.07
$ws.Range("M136").Value = -380093.07

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 33334534
$ws.Range("I81").Value = 55556224
$ws.Range("J81").Value = 1997.5
$ws.Range("K81").Value = 111112448
$ws.Range("L81").Value = 3995
$ws.Range("M81").Value = -111111387
$ws.Range("N81").Value = -6117
$ws.Range("H84").Value = 33334534
$ws.Range("I84").Value = 55556224
$ws.Range("J84").Value = 1997.5
$ws.Range("K84").Value = 555562240
$ws.Range("L84").Value = 19975
$ws.Range("M84").Value = -555556936
$ws.Range("N84").Value = -30583
$ws.Range("H113").Value = 682.9259
$ws.Range("I113").Value = 582.8
$ws.Range("K113").Value = 1748.4
$ws.Range("M113").Value = 421.6000000000001
$ws.Range("H135").Value = 59573.383
$ws.Range("I135").Value = 47000
$ws.Range("J135").Value = 70350.57000000001
$ws.Range("K135").Value = 47000
$ws.Range("L135").Value = 70350.57000000001
$ws.Range("M135").Value = -41930
$ws.Range("N135").Value = -80490.57000000001
$ws.Range("H136").Value = 3463.375
$ws.Range("I136").Value = 1679.2307
$ws.Range("K136").Value = 5037.6921
$ws.Range("M136").Value = -2487.6921
